# Re-pull update for ENAR-Risk-Factors-Analyses/data/FinalVars.xlsx
#
# 1. QuesVars (sheet1): a new questionnaire variable row was pulled in —
#    "ALQ130 - Avg # alcohol drinks/day - past 12 mos" — inserted right
#    after the existing ALQ120Q row (old row 34), pushing every row below
#    it down by one.
# 2. The active window/tab moved from ExamVars back to QuesVars, with a
#    fresh scroll position + cell selection on QuesVars, and ExamVars is
#    no longer the selected tab.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("QuesVars")

# Insert a new row above the old row 34 (ALQ120Q) and populate it with the
# newly-pulled variable. Everything that was at row 34+ shifts down to 35+.
$ws.Rows.Item(34).EntireRow.Insert()
$ws.Cells.Item(34, 1).Value = "ALQ130 - Avg # alcohol drinks/day - past 12 mos"

# Make QuesVars the active sheet/tab again (it had drifted to ExamVars),
# and restore the scrolled/selected view state that was captured on save.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E39").Select()
